# Auto-generated edit script for Guayaba workbook weekly update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C7").Value = "Arica y Parinacota"
$ws.Range("D7").Value = 44425
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100108
$ws.Range("H7").Value = "Tropicales y subtropicales"
$ws.Range("I7").Value = 100108001
$ws.Range("J7").Value = "Guayaba"
$ws.Range("K7").Value = "Sin especificar"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 140
$ws.Range("N7").Value = 1200
$ws.Range("O7").Value = 1300
$ws.Range("P7").Value = 1250
$ws.Range("Q7").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R7").Value = "Región de Arica y Parinacota"
$ws.Range("S7").Value = 1250
$ws.Range("T7").Value = 1

# Row 8
$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C8").Value = "Arica y Parinacota"
$ws.Range("D8").Value = 44407
$ws.Range("E8").Value = 15
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100108
$ws.Range("H8").Value = "Tropicales y subtropicales"
$ws.Range("I8").Value = 100108001
$ws.Range("J8").Value = "Guayaba"
$ws.Range("K8").Value = "Sin especificar"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 600
$ws.Range("O8").Value = 650
$ws.Range("P8").Value = 625
$ws.Range("Q8").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R8").Value = "Región de Arica y Parinacota"
$ws.Range("S8").Value = 625
$ws.Range("T8").Value = 1

# Row 9
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C9").Value = "Arica y Parinacota"
$ws.Range("D9").Value = 44414
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100108
$ws.Range("H9").Value = "Tropicales y subtropicales"
$ws.Range("I9").Value = 100108001
$ws.Range("J9").Value = "Guayaba"
$ws.Range("K9").Value = "Sin especificar"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 160
$ws.Range("N9").Value = 1300
$ws.Range("O9").Value = 1400
$ws.Range("P9").Value = 1350
$ws.Range("Q9").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R9").Value = "Región de Arica y Parinacota"
$ws.Range("S9").Value = 1350
$ws.Range("T9").Value = 1

# Row 10
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C10").Value = "Arica y Parinacota"
$ws.Range("D10").Value = 44403
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100108
$ws.Range("H10").Value = "Tropicales y subtropicales"
$ws.Range("I10").Value = 100108001
$ws.Range("J10").Value = "Guayaba"
$ws.Range("K10").Value = "Sin especificar"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 100
$ws.Range("N10").Value = 1200
$ws.Range("O10").Value = 1300
$ws.Range("P10").Value = 1250
$ws.Range("Q10").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R10").Value = "Región de Arica y Parinacota"
$ws.Range("S10").Value = 1250
$ws.Range("T10").Value = 1

# Row 11
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C11").Value = "Arica y Parinacota"
$ws.Range("D11").Value = 44403
$ws.Range("E11").Value = 15
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100108
$ws.Range("H11").Value = "Tropicales y subtropicales"
$ws.Range("I11").Value = 100108001
$ws.Range("J11").Value = "Guayaba"
$ws.Range("K11").Value = "Sin especificar"
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 120
$ws.Range("N11").Value = 950
$ws.Range("O11").Value = 1000
$ws.Range("P11").Value = 975
$ws.Range("Q11").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R11").Value = "Región de Arica y Parinacota"
$ws.Range("S11").Value = 975
$ws.Range("T11").Value = 1

# Row 12
$ws.Range("A12").Value = 1
$ws.Range("B12").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C12").Value = "Arica y Parinacota"
$ws.Range("D12").Value = 44379
$ws.Range("E12").Value = 15
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100108
$ws.Range("H12").Value = "Tropicales y subtropicales"
$ws.Range("I12").Value = 100108001
$ws.Range("J12").Value = "Guayaba"
$ws.Range("K12").Value = "Sin especificar"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 150
$ws.Range("N12").Value = 700
$ws.Range("O12").Value = 800
$ws.Range("P12").Value = 747
$ws.Range("Q12").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R12").Value = "Región de Arica y Parinacota"
$ws.Range("S12").Value = 747
$ws.Range("T12").Value = 1

# Row 13
$ws.Range("A13").Value = 1
$ws.Range("B13").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C13").Value = "Arica y Parinacota"
$ws.Range("D13").Value = 44379
$ws.Range("E13").Value = 15
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100108
$ws.Range("H13").Value = "Tropicales y subtropicales"
$ws.Range("I13").Value = 100108001
$ws.Range("J13").Value = "Guayaba"
$ws.Range("K13").Value = "Sin especificar"
$ws.Range("L13").Value = "Segunda"
$ws.Range("M13").Value = 140
$ws.Range("N13").Value = 500
$ws.Range("O13").Value = 600
$ws.Range("P13").Value = 543
$ws.Range("Q13").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R13").Value = "Región de Arica y Parinacota"
$ws.Range("S13").Value = 543
$ws.Range("T13").Value = 1

# Row 14
$ws.Range("A14").Value = 1
$ws.Range("B14").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C14").Value = "Arica y Parinacota"
$ws.Range("D14").Value = 44330
$ws.Range("E14").Value = 15
$ws.Range("F14").Value = "Fruta"
$ws.Range("G14").Value = 100108
$ws.Range("H14").Value = "Tropicales y subtropicales"
$ws.Range("I14").Value = 100108001
$ws.Range("J14").Value = "Guayaba"
$ws.Range("K14").Value = "Sin especificar"
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = 1200
$ws.Range("O14").Value = 1300
$ws.Range("P14").Value = 1250
$ws.Range("Q14").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R14").Value = "Región de Arica y Parinacota"
$ws.Range("S14").Value = 1250
$ws.Range("T14").Value = 1

# Row 15
$ws.Range("A15").Value = 1
$ws.Range("B15").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C15").Value = "Arica y Parinacota"
$ws.Range("D15").Value = 44330
$ws.Range("E15").Value = 15
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100108
$ws.Range("H15").Value = "Tropicales y subtropicales"
$ws.Range("I15").Value = 100108001
$ws.Range("J15").Value = "Guayaba"
$ws.Range("K15").Value = "Sin especificar"
$ws.Range("L15").Value = "Segunda"
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 1000
$ws.Range("O15").Value = 1100
$ws.Range("P15").Value = 1050
$ws.Range("Q15").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R15").Value = "Región de Arica y Parinacota"
$ws.Range("S15").Value = 1050
$ws.Range("T15").Value = 1

# Row 16
$ws.Range("A16").Value = 1
$ws.Range("B16").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C16").Value = "Arica y Parinacota"
$ws.Range("D16").Value = 44351
$ws.Range("E16").Value = 15
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100108
$ws.Range("H16").Value = "Tropicales y subtropicales"
$ws.Range("I16").Value = 100108001
$ws.Range("J16").Value = "Guayaba"
$ws.Range("K16").Value = "Sin especificar"
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 700
$ws.Range("O16").Value = 800
$ws.Range("P16").Value = 750
$ws.Range("Q16").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R16").Value = "Región de Arica y Parinacota"
$ws.Range("S16").Value = 750
$ws.Range("T16").Value = 1

# Row 17
$ws.Range("A17").Value = 1
$ws.Range("B17").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C17").Value = "Arica y Parinacota"
$ws.Range("D17").Value = 44351
$ws.Range("E17").Value = 15
$ws.Range("F17").Value = "Fruta"
$ws.Range("G17").Value = 100108
$ws.Range("H17").Value = "Tropicales y subtropicales"
$ws.Range("I17").Value = 100108001
$ws.Range("J17").Value = "Guayaba"
$ws.Range("K17").Value = "Sin especificar"
$ws.Range("L17").Value = "Segunda"
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 600
$ws.Range("O17").Value = 700
$ws.Range("P17").Value = 650
$ws.Range("Q17").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R17").Value = "Región de Arica y Parinacota"
$ws.Range("S17").Value = 650
$ws.Range("T17").Value = 1

# Row 18
$ws.Range("A18").Value = 1
$ws.Range("B18").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C18").Value = "Arica y Parinacota"
$ws.Range("D18").Value = 44350
$ws.Range("E18").Value = 15
$ws.Range("F18").Value = "Fruta"
$ws.Range("G18").Value = 100108
$ws.Range("H18").Value = "Tropicales y subtropicales"
$ws.Range("I18").Value = 100108001
$ws.Range("J18").Value = "Guayaba"
$ws.Range("K18").Value = "Sin especificar"
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 140
$ws.Range("N18").Value = 750
$ws.Range("O18").Value = 800
$ws.Range("P18").Value = 775
$ws.Range("Q18").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R18").Value = "Región de Arica y Parinacota"
$ws.Range("S18").Value = 775
$ws.Range("T18").Value = 1

# Row 19
$ws.Range("A19").Value = 1
$ws.Range("B19").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C19").Value = "Arica y Parinacota"
$ws.Range("D19").Value = 44358
$ws.Range("E19").Value = 15
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100108
$ws.Range("H19").Value = "Tropicales y subtropicales"
$ws.Range("I19").Value = 100108001
$ws.Range("J19").Value = "Guayaba"
$ws.Range("K19").Value = "Sin especificar"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 200
$ws.Range("N19").Value = 700
$ws.Range("O19").Value = 800
$ws.Range("P19").Value = 750
$ws.Range("Q19").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R19").Value = "Región de Arica y Parinacota"
$ws.Range("S19").Value = 750
$ws.Range("T19").Value = 1

# Row 20
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C20").Value = "Arica y Parinacota"
$ws.Range("D20").Value = 44358
$ws.Range("E20").Value = 15
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100108
$ws.Range("H20").Value = "Tropicales y subtropicales"
$ws.Range("I20").Value = 100108001
$ws.Range("J20").Value = "Guayaba"
$ws.Range("K20").Value = "Sin especificar"
$ws.Range("L20").Value = "Segunda"
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = 600
$ws.Range("O20").Value = 650
$ws.Range("P20").Value = 625
$ws.Range("Q20").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R20").Value = "Región de Arica y Parinacota"
$ws.Range("S20").Value = 625
$ws.Range("T20").Value = 1

# Row 21
$ws.Range("A21").Value = 1
$ws.Range("B21").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C21").Value = "Arica y Parinacota"
$ws.Range("D21").Value = 44389
$ws.Range("E21").Value = 15
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100108
$ws.Range("H21").Value = "Tropicales y subtropicales"
$ws.Range("I21").Value = 100108001
$ws.Range("J21").Value = "Guayaba"
$ws.Range("K21").Value = "Sin especificar"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 140
$ws.Range("N21").Value = 750
$ws.Range("O21").Value = 800
$ws.Range("P21").Value = 775
$ws.Range("Q21").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R21").Value = "Región de Arica y Parinacota"
$ws.Range("S21").Value = 775
$ws.Range("T21").Value = 1

# Row 22
$ws.Range("A22").Value = 1
$ws.Range("B22").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C22").Value = "Arica y Parinacota"
$ws.Range("D22").Value = 44389
$ws.Range("E22").Value = 15
$ws.Range("F22").Value = "Fruta"
$ws.Range("G22").Value = 100108
$ws.Range("H22").Value = "Tropicales y subtropicales"
$ws.Range("I22").Value = 100108001
$ws.Range("J22").Value = "Guayaba"
$ws.Range("K22").Value = "Sin especificar"
$ws.Range("L22").Value = "Segunda"
$ws.Range("M22").Value = 120
$ws.Range("N22").Value = 600
$ws.Range("O22").Value = 700
$ws.Range("P22").Value = 650
$ws.Range("Q22").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R22").Value = "Región de Arica y Parinacota"
$ws.Range("S22").Value = 650
$ws.Range("T22").Value = 1

# Row 23
$ws.Range("A23").Value = 1
$ws.Range("B23").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C23").Value = "Arica y Parinacota"
$ws.Range("D23").Value = 44417
$ws.Range("E23").Value = 15
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100108
$ws.Range("H23").Value = "Tropicales y subtropicales"
$ws.Range("I23").Value = 100108001
$ws.Range("J23").Value = "Guayaba"
$ws.Range("K23").Value = "Sin especificar"
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 200
$ws.Range("N23").Value = 1300
$ws.Range("O23").Value = 1400
$ws.Range("P23").Value = 1350
$ws.Range("Q23").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R23").Value = "Región de Arica y Parinacota"
$ws.Range("S23").Value = 1350
$ws.Range("T23").Value = 1

# Row 24
$ws.Range("A24").Value = 1
$ws.Range("B24").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C24").Value = "Arica y Parinacota"
$ws.Range("D24").Value = 44309
$ws.Range("E24").Value = 15
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100108
$ws.Range("H24").Value = "Tropicales y subtropicales"
$ws.Range("I24").Value = 100108001
$ws.Range("J24").Value = "Guayaba"
$ws.Range("K24").Value = "Sin especificar"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 160
$ws.Range("N24").Value = 1400
$ws.Range("O24").Value = 1500
$ws.Range("P24").Value = 1450
$ws.Range("Q24").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R24").Value = "Región de Arica y Parinacota"
$ws.Range("S24").Value = 1450
$ws.Range("T24").Value = 1

# Row 25
$ws.Range("A25").Value = 1
$ws.Range("B25").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C25").Value = "Arica y Parinacota"
$ws.Range("D25").Value = 44344
$ws.Range("E25").Value = 15
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100108
$ws.Range("H25").Value = "Tropicales y subtropicales"
$ws.Range("I25").Value = 100108001
$ws.Range("J25").Value = "Guayaba"
$ws.Range("K25").Value = "Sin especificar"
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 140
$ws.Range("N25").Value = 1000
$ws.Range("O25").Value = 1200
$ws.Range("P25").Value = 1100
$ws.Range("Q25").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R25").Value = "Región de Arica y Parinacota"
$ws.Range("S25").Value = 1100
$ws.Range("T25").Value = 1

# Row 26
$ws.Range("A26").Value = 1
$ws.Range("B26").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C26").Value = "Arica y Parinacota"
$ws.Range("D26").Value = 44344
$ws.Range("E26").Value = 15
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100108
$ws.Range("H26").Value = "Tropicales y subtropicales"
$ws.Range("I26").Value = 100108001
$ws.Range("J26").Value = "Guayaba"
$ws.Range("K26").Value = "Sin especificar"
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 120
$ws.Range("N26").Value = 800
$ws.Range("O26").Value = 850
$ws.Range("P26").Value = 825
$ws.Range("Q26").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R26").Value = "Región de Arica y Parinacota"
$ws.Range("S26").Value = 825
$ws.Range("T26").Value = 1

# Row 27
$ws.Range("A27").Value = 1
$ws.Range("B27").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C27").Value = "Arica y Parinacota"
$ws.Range("D27").Value = 44316
$ws.Range("E27").Value = 15
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100108
$ws.Range("H27").Value = "Tropicales y subtropicales"
$ws.Range("I27").Value = 100108001
$ws.Range("J27").Value = "Guayaba"
$ws.Range("K27").Value = "Sin especificar"
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 140
$ws.Range("N27").Value = 1100
$ws.Range("O27").Value = 1200
$ws.Range("P27").Value = 1150
$ws.Range("Q27").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R27").Value = "Región de Arica y Parinacota"
$ws.Range("S27").Value = 1150
$ws.Range("T27").Value = 1
$ws.Range("D27").NumberFormat = $ws.Range("D2").NumberFormat

# Row 28
$ws.Range("A28").Value = 1
$ws.Range("B28").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C28").Value = "Arica y Parinacota"
$ws.Range("D28").Value = 44326
$ws.Range("E28").Value = 15
$ws.Range("F28").Value = "Fruta"
$ws.Range("G28").Value = 100108
$ws.Range("H28").Value = "Tropicales y subtropicales"
$ws.Range("I28").Value = 100108001
$ws.Range("J28").Value = "Guayaba"
$ws.Range("K28").Value = "Sin especificar"
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 160
$ws.Range("N28").Value = 600
$ws.Range("O28").Value = 700
$ws.Range("P28").Value = 650
$ws.Range("Q28").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R28").Value = "Región de Arica y Parinacota"
$ws.Range("S28").Value = 650
$ws.Range("T28").Value = 1
$ws.Range("D28").NumberFormat = $ws.Range("D2").NumberFormat
